# Fix the calculation issue of the Milestone class.
# Rows 3, 4 and 6 (columns B:O) were computed against the wrong source row.
# Rotate the row contents so that:
#   row 3 <- old row 6
#   row 4 <- old row 3
#   row 6 <- old row 4
# Columns A, P, Q and R are unaffected (identical across these rows already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")

# Capture the current (pre-edit) values for columns B:O of rows 3, 4 and 6
# before any of them get overwritten.
$row3 = @{}
$row4 = @{}
$row6 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("$col`3").Value2
    $row4[$col] = $ws.Range("$col`4").Value2
    $row6[$col] = $ws.Range("$col`6").Value2
}

function Set-CellText($range, $value) {
    # "True"/"False" text is auto-detected as a boolean by the engine (the
    # Value2 getter even hands it back as a .NET [bool]); force it back to
    # plain text so it round-trips the way the workbook originally stored
    # it (inlineStr "True"/"False"), not a boolean cell.
    if ($value -is [bool]) {
        if ($value) { $text = "True" } else { $text = "False" }
        $range.Value2 = "'" + $text
    } elseif ($value -eq "True" -or $value -eq "False") {
        $range.Value2 = "'" + $value
    } else {
        $range.Value2 = $value
    }
}

# Apply the rotation.
foreach ($col in $cols) {
    Set-CellText $ws.Range("$col`3") $row6[$col]
    Set-CellText $ws.Range("$col`4") $row3[$col]
    Set-CellText $ws.Range("$col`6") $row4[$col]
}
